$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.451.15"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.697.85"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("D5").Value = "218.85"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "0.5481"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "0.06450"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "0.07710"
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("D12").Value = "1.694.27"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "4.553"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "0.000008422"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "65.68"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "26.494.74"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "4.948"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D21").Value = "191.26"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "6.261"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "149.06"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  +6.49%  "
$ws.Range("D26").Value = "7.899"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "15.85"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "0.06275"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("D29").Value = "1.379"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "1.331"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "3.615"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "0.6179"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "2.769"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").Value = "0.01646"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").Value = "1.119.01"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "6.123"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "0.8787"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D43").Value = "101.22"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").Value = "1.850.51"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "57.60"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "8.227"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "6.137"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").Value = "0.4300"
$ws.Range("E51").Value = "  -0.09%  "
